# Generate Report for Handoff
# Updates the "Status" text from "In Translation" to "Ready for handoff" and
# refreshes the associated timestamps on the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-08-12 08:46:52"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-08-12 08:46:45"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2016-08-12 08:46:52"

# Widen the columns that contain the newly-lengthened status text, matching
# the column-width growth recorded for the handoff report (status column
# grew from "In Translation" to the longer "Ready for handoff").
$newStatusColumnWidth = 16.25
$wsOverview.Range("E1").EntireColumn.ColumnWidth = $newStatusColumnWidth
$wsOverview.Range("F1").EntireColumn.ColumnWidth = $newStatusColumnWidth
$wsZhCn.Range("C1").EntireColumn.ColumnWidth = $newStatusColumnWidth
$wsDeDe.Range("C1").EntireColumn.ColumnWidth = $newStatusColumnWidth
